$d = $word.ActiveDocument

function Insert-ParaBeforeIndex($idx, $innerXml) {
    $target = $d.Paragraphs.Item($idx)
    $target.Range.InsertParagraphBefore() | Out-Null
    $newPara = $d.Paragraphs.Item($idx)
    $full = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newPara.Range.InsertXML($full) | Out-Null
}

function Insert-ParaAfterIndex($idx, $innerXml) {
    $target = $d.Paragraphs.Item($idx)
    $target.Range.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Item($idx + 1)
    $full = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newPara.Range.InsertXML($full) | Out-Null
}

# Locate the bookmark paragraph (last paragraph in the body, holds the _GoBack bookmark).
$bookmarkIdx = $d.Paragraphs.Count

# --- Insert the first six Q&A pairs before the bookmark paragraph ---
Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Big Data Impacts and Challenges (2019)</w:t></w:r></w:p>'
$bookmarkIdx = $bookmarkIdx + 1
Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:r><w:tab/><w:t xml:space="preserve">What is the problem with the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>data</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$bookmarkIdx = $bookmarkIdx + 1

Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Data Management Challenges for Deep Learning (2019)</w:t></w:r></w:p>'
$bookmarkIdx = $bookmarkIdx + 1
Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:r><w:tab/><w:t xml:space="preserve">How are the limits impacting deep </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>learning</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$bookmarkIdx = $bookmarkIdx + 1

Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Orchestrating Big Data Analysis in the Cloud (2019)</w:t></w:r></w:p>'
$bookmarkIdx = $bookmarkIdx + 1
Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:r><w:tab/><w:t xml:space="preserve">What is the state of the art approaches and how do they </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>differ</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$bookmarkIdx = $bookmarkIdx + 1

Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Patterns for Distributed Real-Time Stream Processing (2017)</w:t></w:r></w:p>'
$bookmarkIdx = $bookmarkIdx + 1
Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:r><w:tab/><w:t xml:space="preserve">How does real-time map/reduce </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>work</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$bookmarkIdx = $bookmarkIdx + 1

Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>The future of FinTech (2019)</w:t></w:r></w:p>'
$bookmarkIdx = $bookmarkIdx + 1
Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:r><w:tab/><w:t>Expand across an industry how these are applicable</w:t></w:r></w:p>'
$bookmarkIdx = $bookmarkIdx + 1

Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>End-to-End Time Architecture for Analyzing and Clustering Timeseries Data (2018)</w:t></w:r></w:p>'
$bookmarkIdx = $bookmarkIdx + 1
Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:r><w:tab/><w:t xml:space="preserve">Provide a concrete example of how it </w:t></w:r><w:r><w:t>is applied</w:t></w:r></w:p>'
$bookmarkIdx = $bookmarkIdx + 1

# --- Insert the heading for the bookmark item (its answer reuses the existing bookmark paragraph) ---
Insert-ParaBeforeIndex $bookmarkIdx '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Supervised Sentiment Analysis of Tweets (2019)</w:t></w:r></w:p>'
$bookmarkIdx = $bookmarkIdx + 1

# --- Fill the bookmark paragraph in place, preserving the _GoBack bookmark ---
$bmPara = $d.Paragraphs.Item($bookmarkIdx)
$bmRange = $bmPara.Range
$bmFull = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t>Extract</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> subsets of information from the firehose.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$bmRange.InsertXML($bmFull) | Out-Null
# InsertXML adds the new content as a fresh paragraph just before the (now empty) original
# bookmark-home paragraph; remove that leftover empty paragraph mark.
$leftoverIdx = $bookmarkIdx + 1
$prevPara = $d.Paragraphs.Item($bookmarkIdx)
$leftoverPara = $d.Paragraphs.Item($leftoverIdx)
$delRange = $d.Range($prevPara.Range.End - 1, $leftoverPara.Range.End)
$delRange.Delete() | Out-Null

# --- Append the remaining Q&A pairs after the bookmark paragraph ---
$afterIdx = $bookmarkIdx
Insert-ParaAfterIndex $afterIdx '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>TICC of Multivariate Time Series Data (2017)</w:t></w:r></w:p>'
$afterIdx = $afterIdx + 1
Insert-ParaAfterIndex $afterIdx '<w:p><w:r><w:tab/><w:t xml:space="preserve">How can we derive context from the time </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>series</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$afterIdx = $afterIdx + 1

Insert-ParaAfterIndex $afterIdx '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Adaptive Deep Learning for Incremental Learning</w:t></w:r></w:p>'
$afterIdx = $afterIdx + 1
Insert-ParaAfterIndex $afterIdx '<w:p><w:r><w:tab/><w:t>Improving deep learning using shorter retentions</w:t></w:r></w:p>'
$afterIdx = $afterIdx + 1

Write-Output "done"
